# Rename the header label in J1 from "e_ps_qs_UC" to "e_qs_ps_UC"
# (part of splitting the combined function into orig/star/hat/bar/tilde variants).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J1").Value2 = "e_qs_ps_UC"

# Update the view so the previously-hidden column D becomes visible again
# and the selection moves from R4 to J1.
$ws.Range("J1").Select()
$aw = $excel.ActiveWindow
$aw.ScrollColumn = 4
$aw.ScrollRow = 1
